$wb = $excel.ActiveWorkbook
$wsPB = $wb.Worksheets.Item("Product Backlog")
$wsInstr = $wb.Worksheets.Item("Instructivo")

# --- Cell content changes on "Product Backlog" ---
$wsPB.Range("C19").Value = "Necesito visualizar historial de clientes que reciclaron en mi negocio"
$wsPB.Range("P14").Value = "17/07/2024 21/07/2026"

# Historias Tecnicas / Historias marked as finished
$wsPB.Range("E24").Value = "Finalizado"
$wsPB.Range("E25").Value = "Finalizado"
$wsPB.Range("E26").Value = "Finalizado"
$wsPB.Range("E27").Value = "Finalizado"
$wsPB.Range("E28").Value = "Finalizado"

# --- Sheet view / selection state ---
$wsPB.Range("C18").Select()
$wsInstr.Activate()

